############################################################
# PM19 Tidsregistrering af Marc.xlsx - update 25/02
#
# Adds three new time-tracking entries (rows 4-6) on "Ark1",
# changes the "Aktuelt tidforbrug" (F) cell on row 3 from a
# numeric value to a descriptive text, fills in F4-F6 with
# descriptive text too, and nudges the saved selection on both
# sheets.
############################################################

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Ark1
$ws2 = $wb.Worksheets.Item(2)   # Ark2

# ------------------------------------------------------------
# Row 3: F3 switches from a plain number (2) to free text.
# ------------------------------------------------------------
$ws1.Range("F3").Value = "2 timer"
$ws1.Range("F3").HorizontalAlignment = -4108   # xlCenter

# ------------------------------------------------------------
# Row 4: new entry - UI-prototype android / reviewer
# ------------------------------------------------------------
$ws1.Range("A4").Value = "UI-prototype android"
$ws1.Range("B4").Value = "reviewer"
$ws1.Range("C4").Value = 43886
$ws1.Range("D4").Value = 0.51388888888888895
$ws1.Range("E4").Value = 0.52083333333333337
$ws1.Range("F4").Value = "10 til 20 min"
$ws1.Range("F4").NumberFormat = "mmm-yy"

# ------------------------------------------------------------
# Row 5: new entry - UI-prototype desktop / reviewer
# ------------------------------------------------------------
$ws1.Range("A5").Value = "UI-prototype desktop"
$ws1.Range("B5").Value = "reviewer"
$ws1.Range("C5").Value = 43886
$ws1.Range("D5").Value = 0.56944444444444442
$ws1.Range("E5").Value = 0.57638888888888895
$ws1.Range("F5").Value = "10 til 20 min"
$ws1.Range("F5").NumberFormat = "mmm-yy"

# ------------------------------------------------------------
# Row 6: new entry - OC0103 - angiv primo / System Analyst
# ------------------------------------------------------------
$ws1.Range("A6").Value = "OC0103 - angiv primo"
$ws1.Range("B6").Value = "System Analyst "
$ws1.Range("B6").Style = "Normal"
$ws1.Range("C6").Value = 43886
$ws1.Range("D6").Value = 0.58333333333333337
$ws1.Range("E6").Value = 0.63194444444444442
$ws1.Range("F6").Value = "1 time"

# ------------------------------------------------------------
# Selection bookkeeping (matches the saved cursor position in
# the authored workbook).
# ------------------------------------------------------------
$ws2.Range("B5").Select()
$ws1.Range("D16").Select()
